# Illustration.pptx / slide 4 ("10 11" mind-map) — widen the legend box,
# extend its caption text, and shift the node/connector layout down-and-right
# to make room (per the authoritative OOXML diff).
#
# PowerPoint's Shape.Left/Top/Width/Height are exposed as single-precision
# (32-bit) floats measured in points, while the underlying OOXML stores EMU
# (1 pt = 12700 EMU) as integers. A naive `emu / 12700.0` round-trips through
# that float32 cast and then gets floored back to EMU, which can land 1 EMU
# below the intended integer. EmuToPt nudges the point value up by tiny
# increments until the float32 -> EMU floor reproduces the exact target EMU,
# so the saved XML matches the diff bit-for-bit instead of being off by 1.
function EmuToPt([double]$emu) {
    $base = $emu / 12700.0
    for ($i = 0; $i -le 300; $i++) {
        $cand = $base + ($i * 0.00001)
        $f = [System.Single]$cand
        $back = [System.Math]::Floor([double]$f * 12700.0)
        if ($back -eq $emu) {
            return $cand
        }
    }
    return $base
}

function SetShapeRect($shp, $offX, $offY, $extCx, $extCy) {
    $shp.Left = EmuToPt $offX
    $shp.Top = EmuToPt $offY
    $shp.Width = EmuToPt $extCx
    $shp.Height = EmuToPt $extCy
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Legend / key rectangle: widen it and append to its first line -------
$legend = $s.Shapes.Item("Rectangle 1")
# Grow the box first (its body is auto-fit) so that appending text to the
# first paragraph doesn't trigger an autofit height re-wrap; only the width
# should grow, height must stay put, matching the diff exactly.
SetShapeRect $legend 1655854 474607 8135417 4524315
$legend.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "10 11 # number of people, relationship between each"

# --- Node ovals + their numeric labels: shift down/right ------------------
SetShapeRect ($s.Shapes.Item("Oval 31")) 4244942 1116029 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 36")) 4437125 1277847 573241 369332
SetShapeRect ($s.Shapes.Item("Oval 37")) 4244942 2521906 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 38")) 4437125 2683724 573241 369332
SetShapeRect ($s.Shapes.Item("Oval 39")) 5667912 2521906 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 40")) 5860095 2683724 573241 369332
SetShapeRect ($s.Shapes.Item("Oval 41")) 2929850 2549742 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 42")) 3125059 2683724 573241 369332
SetShapeRect ($s.Shapes.Item("Oval 43")) 5667912 3921933 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 44")) 5860095 4083751 573241 369332
SetShapeRect ($s.Shapes.Item("Oval 45")) 4244942 3921933 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 46")) 4437125 4060903 573241 369332
SetShapeRect ($s.Shapes.Item("Oval 47")) 6964167 3921933 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 48")) 7156350 4083751 573241 369332
SetShapeRect ($s.Shapes.Item("Oval 51")) 7453900 2521906 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 52")) 7646083 2683724 573241 369332
SetShapeRect ($s.Shapes.Item("Oval 53")) 8876870 2521906 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 54")) 8998243 2683724 573241 369332
SetShapeRect ($s.Shapes.Item("Oval 55")) 8219324 1116029 657546 647272
SetShapeRect ($s.Shapes.Item("TextBox 56")) 8411507 1277847 573241 369332

# --- Straight connectors between the ovals: follow their endpoints --------
SetShapeRect ($s.Shapes.Item("Straight Connector 73")) 3491101 1668510 850136 976023
SetShapeRect ($s.Shapes.Item("Straight Connector 74")) 4806193 1668510 958014 948187
SetShapeRect ($s.Shapes.Item("Straight Connector 75")) 4573715 1763301 0 758605
SetShapeRect ($s.Shapes.Item("Straight Connector 76")) 4573715 3169178 0 752755
SetShapeRect ($s.Shapes.Item("Straight Connector 77")) 4806193 3074387 958014 942337
SetShapeRect ($s.Shapes.Item("Straight Connector 78")) 5996685 3169178 0 752755
SetShapeRect ($s.Shapes.Item("Straight Connector 79")) 4902488 4245569 765424 0

# Connector 80 is re-routed: it used to drop straight down (flipV) between
# connection sites 0/4; now it runs horizontally (flipH) between sites 2/6.
$conn80 = $s.Shapes.Item("Straight Connector 80")
$conn80.ConnectorFormat.BeginConnectionSite = 2
$conn80.ConnectorFormat.EndConnectionSite = 6
$conn80.VerticalFlip = $false
$conn80.HorizontalFlip = $true
SetShapeRect $conn80 6325458 4245569 638709 0

SetShapeRect ($s.Shapes.Item("Straight Connector 81")) 7782673 1668510 532946 853396
SetShapeRect ($s.Shapes.Item("Straight Connector 82")) 8780575 1668510 425068 853396
